# Add new category "on-09-09" / "Anti-CTLA-4" as a new row just above the
# existing "on-99-00" / "Other Antineoplastic Agents" row (row 243), pushing
# all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 243 (shifts rows 243:288 down to 244:289),
# inheriting formatting from the row above as Excel normally does.
$ws.Rows("243:243").Insert()

# Populate the new row with the new category code / name pair.
$ws.Range("A243").Value = "on-09-09"
$ws.Range("B243").Value = "Anti-CTLA-4"

# Update the active selection to mirror where the author last worked.
$ws.Range("E241").Select() | Out-Null
